$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Section" column (header in J1, value "第五事业部" in J2) is being
# removed entirely; Leader/DOB/DOJ (previously K:M) shift left into J:L.
$ws.Range("J1:J2").EntireColumn.Delete()

# Match the author's resulting selection (whole-column selection on the
# new J column, which now holds "Leader").
$ws.Range("J1:J1048576").Select()
